$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.496.40"
$ws.Range("E2").Value = "  +2.93%  "

$ws.Range("D3").Value = "1.606.14"
$ws.Range("E3").Value = "  +2.64%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.57"
$ws.Range("E5").Value = "  +1.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.520"
$ws.Range("E6").Value = "  +6.87%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "26.59"
$ws.Range("E8").Value = "  +5.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.35"
$ws.Range("E9").Value = "  -4.71%  "

$ws.Range("E10").Value = "  +2.48%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0599"
$ws.Range("E11").Value = "  +2.46%  "

$ws.Range("E12").Value = "  +1.92%  "

$ws.Range("D13").Value = "1.834.36"
$ws.Range("E13").Value = "  +2.48%  "

$ws.Range("D14").Value = "1.596.88"
$ws.Range("E14").Value = "  +1.30%  "

$ws.Range("D15").Value = "29.520.64"
$ws.Range("E15").Value = "  +2.87%  "

$ws.Range("E16").Value = "  +3.80%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.69"
$ws.Range("E17").Value = "  +1.47%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.37"
$ws.Range("E18").Value = "  +3.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.67"
$ws.Range("E19").Value = "  +4.72%  "

$ws.Range("E20").Value = "  +3.68%  "

$ws.Range("E21").Value = "  +2.45%  "

$ws.Range("E22").Value = "  +0.25%  "

$ws.Range("E23").Value = "  +1.86%  "

$ws.Range("E24").Value = "  +1.49%  "

$ws.Range("E25").Value = "  -0.63%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.64"
$ws.Range("E26").Value = "  +2.60%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.33"
$ws.Range("E27").Value = "  +3.56%  "

$ws.Range("E28").Value = "  +5.39%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.38"
$ws.Range("E29").Value = "  +2.42%  "

$ws.Range("E30").Value = "  -0.08%  "

$ws.Range("E31").Value = "  +2.75%  "

$ws.Range("E32").Value = "  +0.02%  "

$ws.Range("E33").Value = "  +1.73%  "

$ws.Range("E34").Value = "  +4.08%  "

$ws.Range("D35").Value = "1.416.78"
$ws.Range("E35").Value = "  +1.86%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.03"
$ws.Range("E36").Value = "  -1.18%  "

$ws.Range("E37").Value = "  +2.90%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.81"
$ws.Range("E38").Value = "  +5.89%  "

$ws.Range("E39").Value = "  +0.13%  "

$ws.Range("E40").Value = "  +2.02%  "

$ws.Range("E41").Value = "  +3.77%  "

$ws.Range("E42").Value = "  +0.47%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "53.50"
$ws.Range("E43").Value = "  +23.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0481"
$ws.Range("E44").Value = "  +4.70%  "

$ws.Range("E45").Value = "  -0.04%  "

$ws.Range("E46").Value = "  +2.29%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "65.64"
$ws.Range("E47").Value = "  +2.38%  "

$ws.Range("E48").Value = "  +0.04%  "

$ws.Range("D49").Value = "1.745.99"
$ws.Range("E49").Value = "  +2.37%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "86.66"
$ws.Range("E50").Value = "  +1.64%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.840"
$ws.Range("E51").Value = "  -3.16%  "
